$d = $word.ActiveDocument

# 1. Update the "Interactive CLI Builder" feature row text
$d.Content.Find.Execute(
    "Interactive CLI Builder | User-friendly command creation | ✅ New in 8.4",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Interactive CLI Builder | Interactive command creation with real-time feedback | ✅ Implemented",
    2
)

# 2. Insert a new paragraph for interactive_builder.py before the "Others" row,
#    and update the "Others" row counts/percentage.
$rng = $d.Content
$rng.Find.Execute("Others | 1,480 | 18% | Supporting modules", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.InsertParagraphBefore()
$rng.InsertBefore("interactive_builder.py | 300 | 4% | Interactive CLI Builder functionality")

$d.Content.Find.Execute(
    "Others | 1,480 | 18% | Supporting modules",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Others | 1,180 | 16% | Supporting modules",
    2
)

# 3. Update the "Interactive Builder" row in the examples table
$d.Content.Find.Execute(
    "Interactive Builder | Interactive CLI creation | interactive_builder | 50+",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Interactive Builder | Interactive CLI creation | interactive_builder() | 50",
    2
)
